# The test-data cell value in column D ("test cell" identifier) had a
# typo - it was missing the hyphen that the other two rows' sibling
# value implies ("SNP07 TR1&TR2 Report" header uses the un-hyphenated
# form, but the actual cell identifier should read "S-NP07"). Fix it.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Adhesion_TR2_All")

$ws.Range("D3").Value = "S-NP07"

# Leave the selection where the author last left it before saving.
$null = $ws.Range("E10").Select()
